$wb = $excel.ActiveWorkbook

# The "optimization_parameters" sheet had a stray leftover row (row 16:
# "Sheet" / 3 / 4) that is no longer needed. Remove it entirely - this
# shifts the following row (the simulation_timepoints row) up to become
# row 16, and naturally drops the orphaned "Sheet" shared string / the
# now-unused integer number-format style from the workbook.
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Delete()

# The edit session ended with "network_weights" as the active sheet and
# cell B2 selected.
$wsWeights = $wb.Worksheets.Item("network_weights")
$wsWeights.Activate()
$wsWeights.Range("B2").Select()
